# Generate Report for Handoff
#
# This records a new handoff pass for four files that were previously
# handed back / in translation:
#   - 157ae323-...md  (row 6)  : status stays "Handback transform failed"
#   - 337b147f-...md  (row 8)  : status moves from "In Translation" to "Ready for handoff"
#   - 937b0ea5-...md  (row 9)  : status stays "Ready for handoff"
#   - 9db5387c-...md  (row 10) : status stays "Ready for handoff"
# Their "Latest Handoff Date"/"Latest Handoff Datetime" values are refreshed
# to reflect the new handoff.

$wb = $excel.ActiveWorkbook

$newHandoffDateOverview = "2016-03-23 22:25:29"
$newHandoffDateTimeZhCn = "2016-03-23 22:25:25"
$newHandoffDateTimeDeDe = "2016-03-23 22:25:29"

# ---- Overview sheet ----
$overview = $wb.Worksheets.Item("Overview")

# Row 6: 157ae323-8aae-46d1-a102-fad05cbd9ac9.md -> refresh Latest Handoff Date
$overview.Range("D6").Value = $newHandoffDateOverview

# Row 8: 337b147f-63c0-4161-bdae-0ff3c740797b.md -> now Ready for handoff
$overview.Range("B8").Value = "Ready for handoff"
$overview.Range("C8").Value = "Ready for handoff"
$overview.Range("D8").Value = $newHandoffDateOverview

# Row 9: 937b0ea5-dc62-437e-af51-6ac0c5eb0467.md -> refresh Latest Handoff Date
$overview.Range("D9").Value = $newHandoffDateOverview

# Row 10: 9db5387c-fd65-4243-b7a8-e20158d56769.md -> refresh Latest Handoff Date
$overview.Range("D10").Value = $newHandoffDateOverview

# ---- zh-cn sheet ----
$zhcn = $wb.Worksheets.Item("zh-cn")

# Row 6: 157ae323-...md -> refresh Latest Handoff Datetime
$zhcn.Range("E6").Value = $newHandoffDateTimeZhCn

# Row 8: 337b147f-...md -> status Ready for handoff + refresh datetime
$zhcn.Range("C8").Value = "Ready for handoff"
$zhcn.Range("E8").Value = $newHandoffDateTimeZhCn

# Row 9: 937b0ea5-...md -> refresh Latest Handoff Datetime
$zhcn.Range("E9").Value = $newHandoffDateTimeZhCn

# Row 10: 9db5387c-...md -> refresh Latest Handoff Datetime
$zhcn.Range("E10").Value = $newHandoffDateTimeZhCn

# ---- de-de sheet ----
$dede = $wb.Worksheets.Item("de-de")

# Row 6: 157ae323-...md -> refresh Latest Handoff Datetime
$dede.Range("E6").Value = $newHandoffDateTimeDeDe

# Row 8: 337b147f-...md -> status Ready for handoff + refresh datetime
$dede.Range("C8").Value = "Ready for handoff"
$dede.Range("E8").Value = $newHandoffDateTimeDeDe

# Row 9: 937b0ea5-...md -> refresh Latest Handoff Datetime
$dede.Range("E9").Value = $newHandoffDateTimeDeDe

# Row 10: 9db5387c-...md -> refresh Latest Handoff Datetime
$dede.Range("E10").Value = $newHandoffDateTimeDeDe
